# Automatic update of files.
# - Column C ("Förändrad") for data rows 2..117 changes from 45190 to 45192.
# - Two new data rows (118, 119) are appended with fresh "Avverkningsanmälan" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed/updated) date column for all existing data rows.
for ($r = 2; $r -le 117; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# The previously-last row now sits mid-table; pin its (already default) row
# height explicitly, matching how Excel re-serialises a no-longer-final row.
$ws.Rows.Item(117).RowHeight = 15

# Row 118: A 44654-2023
$ws.Range("A118").Value = "A 44654-2023"
$ws.Range("B118").Value = 45189
$ws.Range("B118").NumberFormat = "YYYY-MM-DD"
$ws.Range("C118").Value = 45192
$ws.Range("C118").NumberFormat = "YYYY-MM-DD"
$ws.Range("D118").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E118").Value = "MULLSJÖ"
$ws.Range("F118").Value = "Kommuner"
$ws.Range("G118").Value = 4.7
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 0
$ws.Range("N118").Value = 0
$ws.Range("O118").Value = 0
$ws.Range("P118").Value = 0
$ws.Range("Q118").Value = 0
$ws.Range("R118").WrapText = $true
$ws.Rows.Item(118).RowHeight = 15

# Row 119: A 44649-2023
$ws.Range("A119").Value = "A 44649-2023"
$ws.Range("B119").Value = 45189
$ws.Range("B119").NumberFormat = "YYYY-MM-DD"
$ws.Range("C119").Value = 45192
$ws.Range("C119").NumberFormat = "YYYY-MM-DD"
$ws.Range("D119").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E119").Value = "MULLSJÖ"
$ws.Range("F119").Value = "Kommuner"
$ws.Range("G119").Value = 0.6
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 0
$ws.Range("N119").Value = 0
$ws.Range("O119").Value = 0
$ws.Range("P119").Value = 0
$ws.Range("Q119").Value = 0
$ws.Range("R119").WrapText = $true
